$wb = $excel.ActiveWorkbook

# Apply scheduled-runner profit/price updates across sheets, per the upstream diff.
# Each hunk below corresponds to one worksheet row whose pricing columns
# (H:currentAveragePrice, I:currentAveragePriceNQ, J:currentAveragePriceHQ,
#  K:LevePriceNQ, L:LevePriceHQ, M:LeveProfitNQ, N:LeveProfitHQ) were refreshed.

$ws = $wb.Worksheets.Item("ALC")
# ALC!row 2
$ws.Range("H2").Value = 833525.5600000001
$ws.Range("I2").Value = 2000192
$ws.Range("J2").Value = 192.42857
$ws.Range("K2").Value = 2000192
$ws.Range("L2").Value = 192.42857
$ws.Range("M2").Value = -2000079
$ws.Range("N2").Value = -418.42857

# ALC!row 9
$ws.Range("H9").Value = 1000236
$ws.Range("I9").Value = 1428881.8
$ws.Range("J9").Value = 62.666668
$ws.Range("K9").Value = 1428881.8
$ws.Range("L9").Value = 62.666668
$ws.Range("M9").Value = -1428712.8
$ws.Range("N9").Value = -400.666668

# ALC!row 58
$ws.Range("H58").Value = 2710.7273
$ws.Range("I58").Value = 165
$ws.Range("J58").Value = 3276.4443
$ws.Range("K58").Value = 495
$ws.Range("L58").Value = 9829.332900000001
$ws.Range("M58").Value = -345
$ws.Range("N58").Value = -10129.3329

# ALC!row 98
$ws.Range("H98").Value = 2284.6155
$ws.Range("I98").Value = 2411.111
$ws.Range("K98").Value = 2411.111
$ws.Range("M98").Value = -913.1109999999999

# ALC!row 103
$ws.Range("H103").Value = 1251030.1
$ws.Range("I103").Value = 2001150.2
$ws.Range("J103").Value = 830
$ws.Range("K103").Value = 6003450.6
$ws.Range("L103").Value = 2490
$ws.Range("M103").Value = -6002864.6
$ws.Range("N103").Value = -3662

# ALC!row 122
$ws.Range("H122").Value = 2284.6155
$ws.Range("I122").Value = 2411.111
$ws.Range("K122").Value = 7233.333
$ws.Range("M122").Value = -4783.333

# ALC!row 129
$ws.Range("H129").Value = 1337.4722
$ws.Range("J129").Value = 1413.9395
$ws.Range("L129").Value = 4241.818499999999
$ws.Range("N129").Value = -14241.8185

# ALC!row 133
$ws.Range("H133").Value = 32930
$ws.Range("J133").Value = 32930
$ws.Range("L133").Value = 32930
$ws.Range("N133").Value = -43050

# ALC!row 136
$ws.Range("H136").Value = 44033.332
$ws.Range("J136").Value = 44033.332
$ws.Range("L136").Value = 44033.332
$ws.Range("N136").Value = -54233.332

# ALC!row 138
$ws.Range("H138").Value = 2829.65
$ws.Range("I138").Value = 869.375
$ws.Range("J138").Value = 3203.0356
$ws.Range("K138").Value = 2608.125
$ws.Range("L138").Value = 9609.106800000001
$ws.Range("M138").Value = 2531.875
$ws.Range("N138").Value = -19889.1068

# ALC!row 140
$ws.Range("H140").Value = 52500
$ws.Range("J140").Value = 52500
$ws.Range("L140").Value = 52500
$ws.Range("N140").Value = -62860

$ws = $wb.Worksheets.Item("ARM")
# ARM!row 34
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

# ARM!row 132
$ws.Range("H132").Value = 2269.6562
$ws.Range("I132").Value = 1879.8889
$ws.Range("J132").Value = 2770.7856
$ws.Range("K132").Value = 5639.6667
$ws.Range("L132").Value = 8312.356800000001
$ws.Range("M132").Value = -3109.6667
$ws.Range("N132").Value = -13372.3568

$ws = $wb.Worksheets.Item("BSM")
# BSM!row 134
$ws.Range("H134").Value = 2067.3684
$ws.Range("I134").Value = 1866.6666
$ws.Range("J134").Value = 2411.4285
$ws.Range("K134").Value = 5599.9998
$ws.Range("L134").Value = 7234.2855
$ws.Range("M134").Value = -3064.9998
$ws.Range("N134").Value = -12304.2855

$ws = $wb.Worksheets.Item("CRP")
# CRP!row 31
$ws.Range("H31").Value = 3375.7222
$ws.Range("I31").Value = 3515.4707
$ws.Range("J31").Value = 1000
$ws.Range("K31").Value = 3515.4707
$ws.Range("L31").Value = 1000
$ws.Range("M31").Value = -3220.4707
$ws.Range("N31").Value = -1590

# CRP!row 34
$ws.Range("H34").Value = 3375.7222
$ws.Range("I34").Value = 3515.4707
$ws.Range("J34").Value = 1000
$ws.Range("K34").Value = 3515.4707
$ws.Range("L34").Value = 1000
$ws.Range("M34").Value = -3313.4707
$ws.Range("N34").Value = -1404

# CRP!row 58
$ws.Range("H58").Value = 7824.1333
$ws.Range("I58").Value = 851.7143
$ws.Range("J58").Value = 13925
$ws.Range("K58").Value = 851.7143
$ws.Range("L58").Value = 13925
$ws.Range("M58").Value = -648.7143
$ws.Range("N58").Value = -14331

# CRP!row 134
$ws.Range("H134").Value = 4100.091
$ws.Range("I134").Value = 3500.3333
$ws.Range("J134").Value = 4325
$ws.Range("K134").Value = 10500.9999
$ws.Range("L134").Value = 12975
$ws.Range("M134").Value = -7965.999899999999
$ws.Range("N134").Value = -18045

# CRP!row 136
$ws.Range("H136").Value = 7824.1333
$ws.Range("I136").Value = 851.7143
$ws.Range("J136").Value = 13925
$ws.Range("K136").Value = 2555.1429
$ws.Range("L136").Value = 41775
$ws.Range("M136").Value = -5.142899999999827
$ws.Range("N136").Value = -46875

$ws = $wb.Worksheets.Item("CUL")
# CUL!row 33
$ws.Range("H33").Value = 160.06667
$ws.Range("I33").Value = 197.36363
$ws.Range("J33").Value = 57.5
$ws.Range("K33").Value = 1184.18178
$ws.Range("L33").Value = 345
$ws.Range("M33").Value = -901.1817799999999
$ws.Range("N33").Value = -911

# CUL!row 117
$ws.Range("H117").Value = 1867.3478
$ws.Range("J117").Value = 1979.55
$ws.Range("L117").Value = 5938.65
$ws.Range("N117").Value = -12822.65

# CUL!row 129
$ws.Range("H129").Value = 1824.8
$ws.Range("I129").Value = 1118
$ws.Range("J129").Value = 2178.2
$ws.Range("K129").Value = 3354
$ws.Range("L129").Value = 6534.599999999999
$ws.Range("M129").Value = 1646
$ws.Range("N129").Value = -16534.6

# CUL!row 131
$ws.Range("H131").Value = 7479217
$ws.Range("I131").Value = 45546704
$ws.Range("J131").Value = 1674.9108
$ws.Range("K131").Value = 136640112
$ws.Range("L131").Value = 5024.732400000001
$ws.Range("M131").Value = -136635072
$ws.Range("N131").Value = -15104.7324

# CUL!row 136
$ws.Range("H136").Value = 2696.8125
$ws.Range("I136").Value = 1649.909
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 4949.727000000001
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = 150.2729999999992
$ws.Range("N136").Value = -25200

$ws = $wb.Worksheets.Item("LTW")
# LTW!row 7
$ws.Range("H7").Value = 2074.6667
$ws.Range("I7").Value = 2046.5
$ws.Range("J7").Value = 2300
$ws.Range("K7").Value = 2046.5
$ws.Range("L7").Value = 2300
$ws.Range("M7").Value = -1934.5
$ws.Range("N7").Value = -2524

# LTW!row 16
$ws.Range("H16").Value = 480
$ws.Range("I16").Value = 480
$ws.Range("K16").Value = 480
$ws.Range("M16").Value = -310

# LTW!row 126
$ws.Range("H126").Value = 2074.6667
$ws.Range("I126").Value = 2046.5
$ws.Range("J126").Value = 2300
$ws.Range("K126").Value = 6139.5
$ws.Range("L126").Value = 6900
$ws.Range("M126").Value = -3669.5
$ws.Range("N126").Value = -11840

# LTW!row 136
$ws.Range("H136").Value = 11000
$ws.Range("I136").Value = 2000
$ws.Range("J136").Value = 14000
$ws.Range("K136").Value = 6000
$ws.Range("L136").Value = 42000
$ws.Range("M136").Value = -3450
$ws.Range("N136").Value = -47100
